$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price / volume data.
# Cells whose new value parses as a plain number (e.g. "39.60") need to be
# force-set as text so Excel keeps the exact text shown in the sheet rather
# than converting to a float and silently dropping trailing zeros; the
# NumberFormat flip + ClearFormats afterwards restores the default (no) style
# so the cell style is left exactly as it was before the edit.

$ws.Range("D2").Value = "49.525.50"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "2.636.11"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "112.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "324.64"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.60"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.87"
$ws.Range("D11").ClearFormats()
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "3.048.34"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("D16").Value = "2.622.77"
$ws.Range("E16").Value = "  -1.64%  "
$ws.Range("E17").Value = "  -3.54%  "
$ws.Range("D18").Value = "49.419.90"
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.95"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "270.47"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.42%  "
$ws.Range("E24").Value = "  -5.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.53"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("E28").Value = "  +3.22%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.70"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.58%  "
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.137"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.46"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("E34").Value = "  +1.99%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.93"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.40%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.04"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.89%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.88"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.11"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.39"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.20"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("E44").Value = "  -4.44%  "
$ws.Range("D45").Value = "2.058.21"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.20"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.11"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +6.07%  "
$ws.Range("E48").Value = "  -4.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.91"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "59.11"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("E51").Value = "  -3.74%  "
